# Applies the "Updated cryptos list" data refresh (Price / Volume(1h) columns,
# plus the #50 coin swap from ONDO to Arweave in row 51) to sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.404.57"
$ws.Range("E2").Value = "'  -0.46%  "
$ws.Range("D3").Value = "'3.462.44"
$ws.Range("E3").Value = "'  -1.59%  "
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'584.26"
$ws.Range("E5").Value = "'  -0.07%  "
$ws.Range("D6").Value = "'177.85"
$ws.Range("E6").Value = "'  +1.15%  "
$ws.Range("E7").Value = "'  +5.54%  "
$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("D9").Value = "'3.462.48"
$ws.Range("E9").Value = "'  -1.48%  "
$ws.Range("D11").Value = "'6.96"
$ws.Range("E11").Value = "'  +0.74%  "
$ws.Range("D12").Value = "'0.418"
$ws.Range("E12").Value = "'  -1.27%  "
$ws.Range("D13").Value = "'4.066.26"
$ws.Range("E13").Value = "'  -1.34%  "
$ws.Range("E14").Value = "'  +1.30%  "
$ws.Range("D15").Value = "'30.19"
$ws.Range("E15").Value = "'  -0.99%  "
$ws.Range("D16").Value = "'66.310.80"
$ws.Range("E16").Value = "'  -0.58%  "
$ws.Range("E17").Value = "'  -0.63%  "
$ws.Range("D18").Value = "'3.478.13"
$ws.Range("E18").Value = "'  -0.95%  "
$ws.Range("D19").Value = "'5.99"
$ws.Range("E19").Value = "'  -1.23%  "
$ws.Range("D20").Value = "'13.87"
$ws.Range("E20").Value = "'  -1.02%  "
$ws.Range("D21").Value = "'371.62"
$ws.Range("E21").Value = "'  -2.76%  "
$ws.Range("D22").Value = "'7.67"
$ws.Range("E22").Value = "'  -3.14%  "
$ws.Range("D23").Value = "'73.28"
$ws.Range("E23").Value = "'  +1.24%  "
$ws.Range("E24").Value = "'  -0.12%  "
$ws.Range("E25").Value = "'  -2.30%  "
$ws.Range("E26").Value = "'  +3.56%  "
$ws.Range("D27").Value = "'10.05"
$ws.Range("E27").Value = "'  +1.62%  "
$ws.Range("D28").Value = "'0.177"
$ws.Range("E28").Value = "'  +2.53%  "
$ws.Range("E29").Value = "'  -0.05%  "
$ws.Range("D30").Value = "'5.98"
$ws.Range("E30").Value = "'  +1.05%  "
$ws.Range("E31").Value = "'  -1.05%  "
$ws.Range("D32").Value = "'23.70"
$ws.Range("E32").Value = "'  -3.71%  "
$ws.Range("E33").Value = "'  -0.08%  "
$ws.Range("E34").Value = "'  -2.69%  "
$ws.Range("E35").Value = "'  -5.77%  "
$ws.Range("E36").Value = "'  -0.90%  "
$ws.Range("D37").Value = "'161.11"
$ws.Range("E37").Value = "'  -0.31%  "
$ws.Range("D38").Value = "'0.886"
$ws.Range("E38").Value = "'  -1.35%  "
$ws.Range("D39").Value = "'28.09"
$ws.Range("E39").Value = "'  -6.82%  "
$ws.Range("E40").Value = "'  +0.98%  "
$ws.Range("D41").Value = "'2.817.79"
$ws.Range("E41").Value = "'  +3.18%  "
$ws.Range("E42").Value = "'  +0.43%  "
$ws.Range("E43").Value = "'  +1.60%  "
$ws.Range("D44").Value = "'6.47"
$ws.Range("E44").Value = "'  -0.74%  "
$ws.Range("D45").Value = "'0.0695"
$ws.Range("E45").Value = "'  -1.17%  "
$ws.Range("D46").Value = "'25.22"
$ws.Range("E46").Value = "'  +0.59%  "
$ws.Range("D47").Value = "'343.01"
$ws.Range("E47").Value = "'  +5.37%  "
$ws.Range("D48").Value = "'40.02"
$ws.Range("E48").Value = "'  -1.76%  "
$ws.Range("E49").Value = "'  -0.48%  "
$ws.Range("E50").Value = "'  +2.83%  "
$ws.Range("B51").Value = "'Arweave"
$ws.Range("C51").Value = "'https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "'31.74"
$ws.Range("E51").Value = "'  +2.39%  "
